$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.392.68'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.148.07'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.25'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.80'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.138.53'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.76'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.28'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000221'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.646.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.308.85'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.146.67'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '513.22'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.80'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.40'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.73'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.09'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.65'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +7.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.83'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.10'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.66'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.60'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '558.48'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.06'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.36'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.32'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0816'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.135.89'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.71%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.25'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.70'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.263'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.97'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.94'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0515'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.08'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.12%  '
